$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "3248"

$t.Cell(6, 1).Range.Text  = "0.36412"
$t.Cell(7, 1).Range.Text  = "0.07607"
$t.Cell(8, 1).Range.Text  = "0.00249"
$t.Cell(9, 1).Range.Text  = "0.35480"
$t.Cell(10, 1).Range.Text = "0.35480"
$t.Cell(11, 1).Range.Text = "0.36412"
$t.Cell(12, 1).Range.Text = "2.19460"

$t.Cell(44, 1).Range.Text = "99.73"
$t.Cell(45, 1).Range.Text = "2.19"
$t.Cell(46, 1).Range.Text = "825"
